$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5197.6
$ws.Range("I40").Value = 4685.8887
$ws.Range("J40").Value = 5965.1665
$ws.Range("K40").Value = 4685.8887
$ws.Range("L40").Value = 5965.1665
$ws.Range("M40").Value = -4510.8887
$ws.Range("N40").Value = -6315.1665

$ws.Range("H99").Value = 779
$ws.Range("I99").Value = 378.75
$ws.Range("K99").Value = 1136.25
$ws.Range("M99").Value = 361.75

$ws.Range("H132").Value = 6337.36
$ws.Range("I132").Value = 7261.3335
$ws.Range("K132").Value = 21784.0005
$ws.Range("M132").Value = -19254.0005

$ws.Range("H137").Value = 2215.5278
$ws.Range("I137").Value = 1763.0454
$ws.Range("J137").Value = 2926.5715
$ws.Range("K137").Value = 5289.1362
$ws.Range("L137").Value = 8779.7145
$ws.Range("M137").Value = -2739.1362
$ws.Range("N137").Value = -13879.7145

$ws.Range("H138").Value = 3185.87
$ws.Range("I138").Value = 3111.0435
$ws.Range("J138").Value = 3208.2207
$ws.Range("K138").Value = 9333.130500000001
$ws.Range("L138").Value = 9624.6621
$ws.Range("M138").Value = -4193.130500000001
$ws.Range("N138").Value = -19904.6621

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1547.39
$ws.Range("I32").Value = 1547.39
$ws.Range("K32").Value = 1547.39
$ws.Range("M32").Value = -1260.39

$ws.Range("H45").Value = 17433.139
$ws.Range("I45").Value = 40400.727
$ws.Range("K45").Value = 40400.727
$ws.Range("M45").Value = -40023.727

$ws.Range("H74").Value = 134927.14
$ws.Range("I74").Value = 208127.11
$ws.Range("K74").Value = 208127.11
$ws.Range("M74").Value = -207253.11

$ws.Range("H77").Value = 134927.14
$ws.Range("I77").Value = 208127.11
$ws.Range("K77").Value = 1040635.55
$ws.Range("M77").Value = -1036267.55

$ws.Range("H132").Value = 2556.5925
$ws.Range("I132").Value = 2115.9744
$ws.Range("K132").Value = 6347.9232
$ws.Range("M132").Value = -3817.9232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7882234
$ws.Range("I105").Value = 403297.2
$ws.Range("K105").Value = 403297.2
$ws.Range("M105").Value = -401550.2

$ws.Range("H108").Value = 90000
$ws.Range("I108").Value = 90000
$ws.Range("K108").Value = 90000
$ws.Range("M108").Value = -86160

$ws.Range("H134").Value = 2542.9707
$ws.Range("I134").Value = 2149.4614
$ws.Range("K134").Value = 6448.3842
$ws.Range("M134").Value = -3913.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2925.4949
$ws.Range("I31").Value = 2730.8416
$ws.Range("J31").Value = 3864.4119
$ws.Range("K31").Value = 2730.8416
$ws.Range("L31").Value = 3864.4119
$ws.Range("M31").Value = -2435.8416
$ws.Range("N31").Value = -4454.4119

$ws.Range("H34").Value = 2925.4949
$ws.Range("I34").Value = 2730.8416
$ws.Range("J34").Value = 3864.4119
$ws.Range("K34").Value = 2730.8416
$ws.Range("L34").Value = 3864.4119
$ws.Range("M34").Value = -2528.8416
$ws.Range("N34").Value = -4268.4119

$ws.Range("H105").Value = 1684.0416
$ws.Range("I105").Value = 935.9375
$ws.Range("J105").Value = 3180.25
$ws.Range("K105").Value = 935.9375
$ws.Range("L105").Value = 3180.25
$ws.Range("M105").Value = 811.0625
$ws.Range("N105").Value = -6674.25

$ws.Range("H122").Value = 3747.625
$ws.Range("I122").Value = 3568.7144
$ws.Range("K122").Value = 10706.1432
$ws.Range("M122").Value = -8256.143199999999

$ws.Range("H134").Value = 2980.9
$ws.Range("I134").Value = 2731.6365
$ws.Range("J134").Value = 3285.5557
$ws.Range("K134").Value = 8194.9095
$ws.Range("L134").Value = 9856.667099999999
$ws.Range("M134").Value = -5659.9095
$ws.Range("N134").Value = -14926.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 944.619
$ws.Range("J107").Value = 946.6667
$ws.Range("L107").Value = 2840.0001
$ws.Range("N107").Value = -6680.0001

$ws.Range("H122").Value = 1241.45
$ws.Range("J122").Value = 1228.1765
$ws.Range("L122").Value = 11053.5885
$ws.Range("N122").Value = -15953.5885

$ws.Range("H129").Value = 2284.92
$ws.Range("J129").Value = 2355.3572
$ws.Range("L129").Value = 7066.071599999999
$ws.Range("N129").Value = -17066.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33473106
$ws.Range("I70").Value = 41839068
$ws.Range("K70").Value = 41839068
$ws.Range("M70").Value = -41838798

$ws.Range("H73").Value = 33473106
$ws.Range("I73").Value = 41839068
$ws.Range("K73").Value = 41839068
$ws.Range("M73").Value = -41838132

$ws.Range("H80").Value = 76925480
$ws.Range("I80").Value = 200002030
$ws.Range("J80").Value = 2625.125
$ws.Range("K80").Value = 200002030
$ws.Range("L80").Value = 2625.125
$ws.Range("M80").Value = -200001032
$ws.Range("N80").Value = -4621.125

$ws.Range("H83").Value = 76925480
$ws.Range("I83").Value = 200002030
$ws.Range("J83").Value = 2625.125
$ws.Range("K83").Value = 1000010150
$ws.Range("L83").Value = 13125.625
$ws.Range("M83").Value = -1000005158
$ws.Range("N83").Value = -23109.625

$ws.Range("H127").Value = 74308.664
$ws.Range("J127").Value = 74308.664
$ws.Range("L127").Value = 74308.664
$ws.Range("N127").Value = -84228.664

$ws.Range("H132").Value = 3578.138
$ws.Range("I132").Value = 3098.8262
$ws.Range("K132").Value = 9296.4786
$ws.Range("M132").Value = -6766.4786

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6784.476
$ws.Range("I7").Value = 4932.4546
$ws.Range("K7").Value = 4932.4546
$ws.Range("M7").Value = -4820.4546

$ws.Range("H82").Value = 5819.643
$ws.Range("I82").Value = 4993.8887
$ws.Range("K82").Value = 4993.8887
$ws.Range("M82").Value = -4632.8887

$ws.Range("H85").Value = 5819.643
$ws.Range("I85").Value = 4993.8887
$ws.Range("K85").Value = 4993.8887
$ws.Range("M85").Value = -3745.8887

$ws.Range("H126").Value = 6784.476
$ws.Range("I126").Value = 4932.4546
$ws.Range("K126").Value = 14797.3638
$ws.Range("M126").Value = -12327.3638

$ws.Range("H136").Value = 7762.5293
$ws.Range("I136").Value = 8434.75
$ws.Range("K136").Value = 25304.25
$ws.Range("M136").Value = -22754.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 782.84
$ws.Range("J113").Value = 805.3
$ws.Range("L113").Value = 2415.9
$ws.Range("N113").Value = -6755.9
